$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

# --- Row 16: previously-empty cells get an explicit "nan" marker, matching
#     the pattern already used by the earlier log rows (13-15) in this sheet.
$ws.Range("B16").Value = "nan"
$ws.Range("C16").Value = "nan"
$ws.Range("D16").Value = "nan"
$ws.Range("E16").Value = "nan"
$ws.Range("F16").Value = "nan"
$ws.Range("G16").Value = "nan"
$ws.Range("H16").Value = "nan"
$ws.Range("I16").Value = "nan"
$ws.Range("J16").Value = "nan"
$ws.Range("K16").Value = "nan"
$ws.Range("M16").Value = "nan"

# --- Row 17: brand new service-log entry for Card21.
# Column A holds "21" as text (matching every other row in this column), so
# build it via a text formula and flatten it to a value to avoid the COM
# layer inferring a numeric type (and to avoid baking a new "Text" number
# format into the cell's style).
$ws.Range("A17").Formula = "=""21"""
$ws.Range("A17").Copy()
$ws.Range("A17").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("L17").Value = "19\12\2024"
$ws.Range("M17").Value = "4320 h"
$ws.Range("N17").Value = "تم عمل  صيانه نصف سنويه"
$ws.Range("O17").Value = "تيم العمل"
